$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reverse the period order in column E (rows 16-24) and swap the
# "Valor Mora" values in F16/F24 to match.
$ws.Range("E16").Value = "2102"
$ws.Range("E17").Value = "2101"
$ws.Range("E18").Value = "2012"
$ws.Range("E19").Value = "2011"
$ws.Range("E20").Value = "2010"
$ws.Range("E21").Value = "2009"
$ws.Range("E22").Value = "2008"
$ws.Range("E23").Value = "2007"
$ws.Range("E24").Value = "2006"

$ws.Range("F16").Value = 25749
$ws.Range("F24").Value = 35112
